# Regenerate the comprehensive_results workbook with a fresh model run
# (new CV folds / hyperparameter search) and add a new "Comparative_Results"
# sheet comparing no-FE / FE / FE+SMOTE scenarios.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Model_Comparison — CV summary stats per model (rows re-scored)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Model_Comparison")
$modelComparison = @(
    @("CatBoost",           0.9408122800714469, 0.004439423084684349, 0.9352297255627505, 0.9474195430465122),
    @("XGBoost",            0.9389822350837246, 0.003817224862165264, 0.9332272773993403, 0.9438348619571052),
    @("GBM",                0.9365267018141573, 0.005620164765551236, 0.9287509694540222, 0.9435562090130458),
    @("LightGBM",           0.9358021757209306, 0.00476923411761053,  0.9291677178818716, 0.9421326051949871),
    @("RandomForest",       0.931586238806392,  0.003900995165645471, 0.9273974714770274, 0.936930150061845),
    @("AdaBoost",           0.9179217360518358, 0.005741320012973289, 0.9092586363169157, 0.9245844710494995),
    @("SVC",                0.8950952007122488, 0.006587526915806162, 0.8866168624262981, 0.9024289948423534),
    @("NeuralNet",          0.8829621914172719, 0.008521545835548957, 0.8709402068791521, 0.892567854559033),
    @("DecisionTree",       0.8679995172648246, 0.005312929410991729, 0.8621174744671507, 0.876990828257369),
    @("LogisticRegression", 0.8578936439703041, 0.008293769420679668, 0.8465987348041002, 0.871335153679199)
)
for ($i = 0; $i -lt $modelComparison.Length; $i++) {
    $row = $modelComparison[$i]
    $r = $i + 2
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
}

# ---------------------------------------------------------------------
# 2. Test_Metrics — single-row test-set metrics for the best model
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Test_Metrics")
$ws2.Range("A2").Value = 0.8577502899110939
$ws2.Range("B2").Value = 0.8560431100846805
$ws2.Range("C2").Value = 0.860015467904099
$ws2.Range("D2").Value = 0.8580246913580247
$ws2.Range("E2").Value = 0.938701855550814

# ---------------------------------------------------------------------
# 3. Feature_Importance — re-ranked feature list with new importances
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Feature_Importance")
$featureImportance = @(
    @("tenure",                                 0.09831535075976734),
    @("InternetService_Fiber optic",             0.08452104515737395),
    @("Contract_Two year",                       0.06505402468366005),
    @("Churn_Risk_Score",                        0.06423153347626288),
    @("PaperlessBilling_Yes",                    0.06207988669831658),
    @("PaymentMethod_Electronic check",          0.05524395818037838),
    @("Contract_One year",                       0.04605512606470746),
    @("TechSupport_Yes",                         0.04476238415211353),
    @("gender_Male",                             0.04283704625333935),
    @("MultipleLines_Yes",                       0.03879958983533133),
    @("MonthlyCharges",                          0.03842125945634847),
    @("OnlineSecurity_Yes",                      0.03838571064059735),
    @("Partner_Yes",                             0.03293259049817319),
    @("StreamingTV_Yes",                         0.03088697596257657),
    @("StreamingMovies_Yes",                     0.03050738294890358),
    @("Engagement_Score",                        0.02984503193220509),
    @("TotalCharges",                            0.02852663208842904),
    @("Dependents_Yes",                          0.02807371620597136),
    @("OnlineBackup_Yes",                        0.02728885202204224),
    @("Service_Utilization",                     0.01912299871012789),
    @("SeniorCitizen",                           0.01859785916434884),
    @("DeviceProtection_Yes",                    0.01670343486910839),
    @("Payment_Reliability",                     0.01641862046388125),
    @("PaymentMethod_Mailed check",              0.01554130097768571),
    @("PaymentMethod_Credit card (automatic)",   0.008123237475500675),
    @("OnlineBackup_No internet service",        0.005798559770640363),
    @("PhoneService_Yes",                        0.003853627837621522),
    @("MultipleLines_No phone service",          0.003160258522460472),
    @("OnlineSecurity_No internet service",      0.003086195260744489),
    @("DeviceProtection_No internet service",    0.002130485645282485),
    @("StreamingTV_No internet service",         0.0006953242861001574),
    @("TechSupport_No internet service",         0),
    @("StreamingMovies_No internet service",     0),
    @("InternetService_No",                      0)
)
for ($i = 0; $i -lt $featureImportance.Length; $i++) {
    $row = $featureImportance[$i]
    $r = $i + 2
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
}

# ---------------------------------------------------------------------
# 4. Hyperparameter_Tuning — updated best score / params for a few models
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Hyperparameter_Tuning")
$ws4.Range("B4").Value  = 0.9351456164279991
$ws4.Range("C4").Value  = "{'n_estimators': 264, 'learning_rate': 0.040435977764864704, 'max_depth': 8, 'subsample': 0.9396810473416601}"
$ws4.Range("B10").Value = 0.9334813234350913
$ws4.Range("C10").Value = "{'n_estimators': 337, 'learning_rate': 0.020713257301171004, 'max_depth': 9}"
$ws4.Range("B11").Value = 0.9361114074765897
$ws4.Range("C11").Value = "{'iterations': 336, 'learning_rate': 0.08491024370403978, 'depth': 10}"

# ---------------------------------------------------------------------
# 5. Summary — headline run summary
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Summary")
$ws5.Range("D2").Value = 0.9408122800714469
$ws5.Range("E2").Value = 0.8577502899110939
$ws5.Range("F2").Value = 0.8560431100846805
$ws5.Range("G2").Value = 0.860015467904099
$ws5.Range("H2").Value = 0.8580246913580247
$ws5.Range("I2").Value = 0.938701855550814
$ws5.Range("J2").Value = "0.31 seconds"

# ---------------------------------------------------------------------
# 6. New sheet — Comparative_Results (no-FE vs FE vs FE+SMOTE)
# ---------------------------------------------------------------------
$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6 = $wb.Worksheets.Add($null, $wsLast)
$ws6.Name = "Comparative_Results"

$ws6.Range("A1").Value = "Scenario"
$ws6.Range("B1").Value = "Accuracy"
$ws6.Range("C1").Value = "Precision"
$ws6.Range("D1").Value = "Recall"
$ws6.Range("E1").Value = "F1"
$ws6.Range("F1").Value = "ROC_AUC"
$ws6.Range("A1:F1").Font.Bold = $true
$ws6.Range("A1:F1").HorizontalAlignment = -4108
$ws6.Range("A1:F1").VerticalAlignment = -4160
$ws6.Range("A1:F1").Borders.LineStyle = 1

$comparative = @(
    @("No Feature Engineering",     0.7847813742191937, 0.6170212765957447, 0.4967880085653105, 0.5504151838671412, 0.8248811017080976),
    @("Feature Engineering",        0.7802385008517888, 0.6104972375690608, 0.4732334047109208, 0.5331724969843185, 0.8239866754482059),
    @("Feature Engineering + SMOTE",0.8577502899110939, 0.8560431100846805, 0.860015467904099,  0.8580246913580247, 0.938701855550814)
)
for ($i = 0; $i -lt $comparative.Length; $i++) {
    $row = $comparative[$i]
    $r = $i + 2
    $ws6.Cells.Item($r, 1).Value = $row[0]
    $ws6.Cells.Item($r, 2).Value = $row[1]
    $ws6.Cells.Item($r, 3).Value = $row[2]
    $ws6.Cells.Item($r, 4).Value = $row[3]
    $ws6.Cells.Item($r, 5).Value = $row[4]
    $ws6.Cells.Item($r, 6).Value = $row[5]
}
